# Update "想去人数" (want-to-go count) figures for two conventions / upcoming
# events, as refreshed by the gh-pages data generation run at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 2442
$ws1.Range("F9").Value = 6404
$ws1.Range("F10").Value = 171

# Sheet "全部类型" (all types) mirrors the same events a couple rows lower
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 2442
$ws4.Range("F11").Value = 6404
$ws4.Range("F12").Value = 171
